# Apply the "data with feed columns" edit to the ILRI GES master template.
#
# Summary of the change (per the OOXML diff):
#  - Header "Feed" (M2) is renamed "Feed Offered"; a new header column N2
#    "Feed Refusal" is added, extending the J1:M1 "Week 1" merge to J1:N1
#    and widening column N to match the other week columns.
#  - The three sample/demo data rows (old rows 3-5) are removed.
#  - Four blank spacer rows are (re)introduced above the hidden
#    "template" row (old row 10, a full A:AK formatted blank row used for
#    appending new week blocks) so that it lands on row 10 again, and the
#    lone formatted cell that used to sit at J11 now sits at J7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Feed Offered / Feed Refusal header columns -------------------------

$ws.Range("M2").Value = "Feed Offered"

# Extend the "Week 1" merged header from J1:M1 to J1:N1.
$ws.Range("J1:M1").UnMerge() | Out-Null
$ws.Range("N1").Value = $ws.Range("M1").Value
$ws.Range("J1:N1").Merge() | Out-Null

# New header cell for the Feed Refusal column.
$ws.Range("N2").Value = "Feed Refusal"
$ws.Range("N2").Font.Bold = $true
$ws.Range("N2").Font.Name = "Calibri"
$ws.Range("N2").Font.Size = 11
$ws.Range("N2").Interior.Color = $ws.Range("M2").Interior.Color
$ws.Range("N2").HorizontalAlignment = $ws.Range("M2").HorizontalAlignment
$ws.Range("N2").VerticalAlignment = $ws.Range("M2").VerticalAlignment
$ws.Range("N2").Borders.LineStyle = $ws.Range("M2").Borders.LineStyle

# Column N should match the width of the other week data columns (col J).
$ws.Columns("N").ColumnWidth = 6.46

# --- 2. Remove the three sample data rows (old rows 3-5) -------------------

$ws.Rows("3:5").Delete() | Out-Null

# --- 3. Re-open spacer rows 6-9 above the template row (old row 10) --------

$ws.Rows("7:9").Insert() | Out-Null

# The formatted-but-empty cell that used to be at J11 now belongs at J7.
$ws.Range("J7").Font.Bold = $true
$ws.Range("J7").Font.Name = "Calibri"
$ws.Range("J7").Font.Size = 11
$ws.Range("J11").ClearContents() | Out-Null
$ws.Range("J11").ClearFormats() | Out-Null

# Re-register the other blank formatted placeholder cells.
$ws.Range("F6").NumberFormat = "General"
$ws.Range("F8").NumberFormat = "General"
$ws.Range("F9:X9").NumberFormat = "General"
$ws.Range("O9:X9").NumberFormat = "General"

# --- 4. Selection / active cell --------------------------------------------

$ws.Range("E7").Select() | Out-Null
